$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain exact text formatting (no numeric auto-conversion)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.492.00"
$ws.Range("E2").Value = "  -1.24%  "

$ws.Range("D3").Value = "1.850.13"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "243.26"
$ws.Range("E5").Value = "  -1.31%  "

$ws.Range("D6").Value = "0.6614"
$ws.Range("E6").Value = "  +4.04%  "

$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "48.04"
$ws.Range("E8").Value = "  +3.11%  "

$ws.Range("D9").Value = "0.07498"

$ws.Range("D10").Value = "0.2993"
$ws.Range("E10").Value = "  -0.12%  "

$ws.Range("D11").Value = "24.40"
$ws.Range("E11").Value = "  -0.34%  "

$ws.Range("D12").Value = "0.07641"
$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").Value = "1.846.56"
$ws.Range("E13").Value = "  -1.07%  "

$ws.Range("D14").Value = "5.023"
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").Value = "0.6855"
$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").Value = "83.83"
$ws.Range("E16").Value = "  -0.60%  "

$ws.Range("D17").Value = "0.000009562"
$ws.Range("E17").Value = "  +1.98%  "

$ws.Range("D18").Value = "6.151"
$ws.Range("E18").Value = "  +0.96%  "

$ws.Range("D19").Value = "29.519.24"
$ws.Range("E19").Value = "  -1.07%  "

$ws.Range("D20").Value = "2.081.06"
$ws.Range("E20").Value = "  -1.54%  "

$ws.Range("D21").Value = "237.07"
$ws.Range("E21").Value = "  -0.78%  "

$ws.Range("E22").Value = "  -0.58%  "

$ws.Range("D24").Value = "7.690"
$ws.Range("E24").Value = "  +4.56%  "

$ws.Range("D25").Value = "0.9998"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").Value = "0.1427"
$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("D27").Value = "156.76"
$ws.Range("E27").Value = "  -1.59%  "

$ws.Range("D28").Value = "8.496"
$ws.Range("E28").Value = "  -0.92%  "

$ws.Range("E29").Value = "  -0.88%  "

$ws.Range("D30").Value = "0.06027"
$ws.Range("E30").Value = "  -1.00%  "

$ws.Range("D31").Value = "1.488"
$ws.Range("E31").Value = "  -1.20%  "

$ws.Range("D32").Value = "1.249"
$ws.Range("E32").Value = "  -1.51%  "

$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("D34").Value = "4.077"
$ws.Range("E34").Value = "  -1.52%  "

$ws.Range("D35").Value = "1.181"
$ws.Range("E35").Value = "  +1.37%  "

$ws.Range("D36").Value = "1.855"
$ws.Range("E36").Value = "  -0.54%  "

$ws.Range("D37").Value = "0.7235"
$ws.Range("E37").Value = "  -0.71%  "

$ws.Range("D38").Value = "2.594"
$ws.Range("E38").Value = "  -0.97%  "

$ws.Range("D39").Value = "2.801"
$ws.Range("E39").Value = "  -1.92%  "

$ws.Range("D40").Value = "0.01780"
$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("D41").Value = "1.199.77"
$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("D42").Value = "6.242"
$ws.Range("E42").Value = "  -1.19%  "

$ws.Range("E43").Value = "  -2.25%  "

$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("D45").Value = "2.010.30"
$ws.Range("E45").Value = "  -0.69%  "

$ws.Range("D46").Value = "101.93"
$ws.Range("E46").Value = "  -0.44%  "

$ws.Range("D47").Value = "66.17"
$ws.Range("E47").Value = "  -0.28%  "

$ws.Range("D48").Value = "7.454"
$ws.Range("E48").Value = "  +11.00%  "

$ws.Range("D49").Value = "0.4061"
$ws.Range("E49").Value = "  -0.63%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.00000000117"
$ws.Range("E50").Value = "  -5.46%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.084"
$ws.Range("E51").Value = "  -2.28%  "
